$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos price/volume table (GitHub Actions data refresh).
# Price values that look like plain numbers are written with a leading
# apostrophe so Excel keeps them as text (matching the original
# inlineStr cells) instead of re-parsing them into floating point
# numbers and losing exact formatting (trailing zeros, etc.).
$ws.Range('D2').Value = '65.475.24'
$ws.Range('E2').Value = '  +1.46%  '
$ws.Range('D3').Value = '3.390.93'
$ws.Range('E3').Value = '  +1.00%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''559.11'
$ws.Range('E5').Value = '  -0.31%  '
$ws.Range('D6').Value = '''175.59'
$ws.Range('E6').Value = '  -0.29%  '
$ws.Range('E7').Value = '  +1.07%  '
$ws.Range('D8').Value = '3.378.89'
$ws.Range('E8').Value = '  +0.91%  '
$ws.Range('E9').Value = '  -0.04%  '
$ws.Range('D10').Value = '''0.170'
$ws.Range('E10').Value = '  +3.58%  '
$ws.Range('D11').Value = '''0.636'
$ws.Range('E11').Value = '  +0.87%  '
$ws.Range('D12').Value = '''53.52'
$ws.Range('E12').Value = '  -2.82%  '
$ws.Range('D13').Value = '''0.0000278'
$ws.Range('E13').Value = '  +0.70%  '
$ws.Range('D14').Value = '''9.21'
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').Value = '3.937.58'
$ws.Range('E15').Value = '  +0.98%  '
$ws.Range('D16').Value = '''18.30'
$ws.Range('E16').Value = '  +0.12%  '
$ws.Range('E17').Value = '  +1.73%  '
$ws.Range('D18').Value = '3.397.98'
$ws.Range('E18').Value = '  +1.22%  '
$ws.Range('D19').Value = '65.572.56'
$ws.Range('E19').Value = '  +1.69%  '
$ws.Range('D20').Value = '''11.84'
$ws.Range('E20').Value = '  +0.16%  '
$ws.Range('D21').Value = '''0.999'
$ws.Range('E21').Value = '  +1.10%  '
$ws.Range('D22').Value = '''464.59'
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').Value = '''4.96'
$ws.Range('E23').Value = '  +2.30%  '
$ws.Range('D24').Value = '''4.12'
$ws.Range('E24').Value = '  +0.23%  '
$ws.Range('D25').Value = '''14.31'
$ws.Range('E25').Value = '  +6.37%  '
$ws.Range('D26').Value = '''87.12'
$ws.Range('E26').Value = '  +1.02%  '
$ws.Range('D27').Value = '''2.92'
$ws.Range('E27').Value = '  +2.69%  '
$ws.Range('D28').Value = '''10.70'
$ws.Range('E28').Value = '  -1.49%  '
$ws.Range('D29').Value = '''8.72'
$ws.Range('E29').Value = '  -0.70%  '
$ws.Range('D30').Value = '''31.06'
$ws.Range('E30').Value = '  +2.98%  '
$ws.Range('D31').Value = '''6.53'
$ws.Range('E31').Value = '  -2.35%  '
$ws.Range('D32').Value = '''63.35'
$ws.Range('E32').Value = '  +7.08%  '
$ws.Range('D33').Value = '''11.48'
$ws.Range('E33').Value = '  -0.08%  '
$ws.Range('D34').Value = '''575.55'
$ws.Range('E34').Value = '  -0.81%  '
$ws.Range('E35').Value = '  -0.23%  '
$ws.Range('E36').Value = '  +0.03%  '
$ws.Range('E37').Value = '  +3.71%  '
$ws.Range('E38').Value = '  +1.07%  '
$ws.Range('D39').Value = '''35.92'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('D40').Value = '''0.373'
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').Value = '0.0₃0740'
$ws.Range('E41').Value = '  -2.12%  '
$ws.Range('D42').Value = '3.124.81'
$ws.Range('E42').Value = '  +1.12%  '
$ws.Range('D43').Value = '''2.80'
$ws.Range('E43').Value = '  -1.17%  '
$ws.Range('D44').Value = '''0.0417'
$ws.Range('E44').Value = '  +0.97%  '
$ws.Range('B45').Value = 'ApeXProtocol'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D45').Value = '''3.19'
$ws.Range('E45').Value = '  -0.50%  '
$ws.Range('B46').Value = 'Stellar'
$ws.Range('C46').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D46').Value = '''0.134'
$ws.Range('E46').Value = '  +1.55%  '
$ws.Range('B47').Value = 'Fetch.AI'
$ws.Range('C47').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D47').Value = '''2.44'
$ws.Range('E47').Value = '  -2.72%  '
$ws.Range('D48').Value = '''1.00'
$ws.Range('E48').Value = '  +0.20%  '
$ws.Range('D49').Value = '''140.13'
$ws.Range('E49').Value = '  +2.98%  '
$ws.Range('D50').Value = '''2.56'
$ws.Range('E50').Value = '  -1.62%  '
$ws.Range('D51').Value = '''8.44'
$ws.Range('E51').Value = '  +1.03%  '
